$wb = $excel.ActiveWorkbook

# Rename "Team_Management" sheet to "Plan"
$planSheet = $wb.Worksheets.Item("Team_Management")
$planSheet.Name = "Plan"

# Move selection within Plan sheet to D13 and activate that sheet/tab
$planSheet.Range("D13").Select()
$planSheet.Activate()

# Move selection on Phase_Dates back to F3 (already there), deselect tab there
$phaseDates = $wb.Worksheets.Item("Phase_Dates")
$phaseDates.Range("F3").Select()

# Re-activate Plan as the final active sheet/tab
$planSheet.Activate()
$planSheet.Range("D13").Select()
